$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.571.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.19%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.439.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.62%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.76%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.537"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.436.59"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.27%  "

$ws.Range("E10").Value = "  +1.70%  "

$ws.Range("E11").Value = "  +1.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.75%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.71%  "

$ws.Range("E15").Value = "  +5.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.881.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.462.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₅0104"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +263.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.436.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.92%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "325.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.48%  "

$ws.Range("E23").Value = "  +1.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.90%  "

$ws.Range("E25").Value = "  -0.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "630.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.77%  "

$ws.Range("E28").Value = "  +13.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0977"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.558.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.13%  "

$ws.Range("E33").Value = "  +5.24%  "

$ws.Range("E34").Value = "  +5.84%  "

$ws.Range("E35").Value = "  +2.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.62%  "

$ws.Range("E37").Value = "  +0.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.11%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.02%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "152.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.25%  "

$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.372"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.24%  "

$ws.Range("E45").Value = "  +1.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("E47").Value = "  +28.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "143.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.62%  "

$ws.Range("E51").Value = "  +1.79%  "
